$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add a new row 34 for the new bug entry, copying the formatting from the
# last existing data row (row 33) so the new row matches the table style.
$ws.Range("A33").Copy()
$ws.Range("A34").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A34").Value = 1

$ws.Range("B33").Copy()
$ws.Range("B34").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B34").Value = "Prevent popup from the webbrowser control"

$excel.CutCopyMode = 0

$ws.Range("B35").Select()
